$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1: scroll the view a bit (selection stays at D9) ---
# (view-only scroll position; selection unchanged)
$ws1.Range("D9").Select()

# --- Sheet2: update existing rows 6 & 7 (block size / hidden size bump) ---
$ws2.Range("B6").Value = 4
$ws2.Range("C6").Value = 512

# Row 7 moves from the EP160 config to the EP130 config (and same block/hidden bump).
# Write A7 first so the new "EP130" shared string is interned before the row
# 8/9 strings, matching the order new strings were appended upstream.
$ws2.Range("A7").Value = "EP130"
$ws2.Range("B7").Value = 4
$ws2.Range("C7").Value = 512

# --- Sheet2: fill in two new experiment rows (8 & 9) ---
# Row 9 first (EP140) then row 8 (EP190) so the new shared strings land in the
# same order as upstream: EP130, EP140, (140,), EP190, (190,).
$ws2.Range("A9").Value = "EP140"
$ws2.Range("B9").Value = 2
$ws2.Range("C9").Value = 512
$ws2.Range("D9").Value = 128
$ws2.Range("E9").Value = 128
$ws2.Range("F9").Value = 0.001
$ws2.Range("G9").Value = 1
$ws2.Range("H9").Value = "(140,)"
$ws2.Range("I9").Value = 140
$ws2.Range("J9").Value = 0
$ws2.Range("B9:G9").HorizontalAlignment = -4131
$ws2.Range("I9").HorizontalAlignment = -4131

$ws2.Range("A8").Value = "EP190"
$ws2.Range("B8").Value = 2
$ws2.Range("C8").Value = 512
$ws2.Range("D8").Value = 128
$ws2.Range("E8").Value = 128
$ws2.Range("F8").Value = 0.01
$ws2.Range("G8").Value = 1
$ws2.Range("H8").Value = "(190,)"
$ws2.Range("I8").Value = 190
$ws2.Range("J8").Value = 0
$ws2.Range("B8:G8").HorizontalAlignment = -4131
$ws2.Range("I8").HorizontalAlignment = -4131

# --- Sheet2: extend used range with a new blank formatted row 25 ---
$ws2.Range("J24").Copy($ws2.Range("J25"))

# --- Sheet2: update the remembered selection ---
$ws2.Range("C11").Select()
